$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '27.949.16'
$ws.Range("D3").Value = '1.768.93'
$ws.Range("E3").Value = '  +0.96%  '
$ws.Range("E4").Value = '  -0.17%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '328.93'
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = '  +1.37%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4563'
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = '  +1.60%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3527'
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = '  -1.03%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '42.14'
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = '  +1.52%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.07393'
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = '  -1.01%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '1.097'
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = '  +1.46%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '1.001'
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = '  -0.17%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '20.74'
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = '  +0.14%  '
$ws.Range("E14").Value = '  +0.52%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '7.197'
$ws.Range("D15").ClearFormats()
$ws.Range("D16").Value = '1.773.00'
$ws.Range("E16").Value = '  +0.88%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '92.66'
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = '  -1.05%  '
$ws.Range("E18").Value = '  +0.52%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.06446'
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = '  +1.05%  '
$ws.Range("E20").Value = '  -0.12%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '16.94'
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = '  -1.18%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '5.777'
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = '  +0.78%  '
$ws.Range("D23").Value = '27.972.39'
$ws.Range("E23").Value = '  +1.22%  '
$ws.Range("E24").Value = '  +0.22%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.148'
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = '  +2.84%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '160.79'
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = '  -3.07%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '20.17'
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = '  +0.19%  '
$ws.Range("D28").Value = '1.976.13'
$ws.Range("E28").Value = '  +0.95%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '124.17'
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = '  -1.04%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.079'
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = '  -1.08%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.09261'
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = '  +0.90%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '5.603'
$ws.Range("D33").ClearFormats()
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '3.662'
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = '  -0.07%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '11.86'
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = '  +1.20%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.02282'
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = '  -0.15%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.06133'
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = '  +1.85%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.2090'
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = '  -0.15%  '
$ws.Range("B39").Value = 'InternetComputer(DFINITY)'
$ws.Range("C39").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '4.937'
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = '  +0.27%  '
$ws.Range("B40").Value = 'TheSandbox'
$ws.Range("C40").Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.6264'
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = '  -0.35%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '1.183'
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = '  +0.21%  '
$ws.Range("E42").Value = '  -0.47%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '7.858'
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = '  +1.28%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '13.23'
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = '  +0.27%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '3.734'
$ws.Range("D45").ClearFormats()
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.5865'
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = '  +0.02%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '122.65'
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = '  +0.33%  '
$ws.Range("E48").Value = '  +0.14%  '
$ws.Range("E49").Value = '  +0.29%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.06822'
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = '  -0.94%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '73.59'
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = '  +2.77%  '
